$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overal Stats")
Write-Host $ws.Name
Write-Host $ws.Range("A1").Value
